$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style / number format to reuse for the "Fecha" (date) column, copied from
# an existing data row (D199) so the new date cells match the workbook's
# existing date formatting.
$dateFormat = $ws.Cells.Item(199, 4).NumberFormat

$rows = @(
    @{ Row=200; A=2; B="Comercializadora del Agro de Limarí"; C="Coquimbo"; D=45239; E=4; F="Fruta"; G=100107; H="Otros"; I=100107002; J="Chirimoya"; K="Cultivar IV Región"; L="Especial"; M=300; N=15000; O=16000; P=15500; Q="$/bandeja 10 kilos"; R="Provincia de Limarí"; S=1550; T=10 },
    @{ Row=201; A=2; B="Comercializadora del Agro de Limarí"; C="Coquimbo"; D=45239; E=4; F="Fruta"; G=100107; H="Otros"; I=100107002; J="Chirimoya"; K="Cultivar IV Región"; L="Primera";  M=360; N=11000; O=12000; P=11500; Q="$/bandeja 10 kilos"; R="Provincia de Limarí"; S=1150; T=10 },
    @{ Row=202; A=2; B="Comercializadora del Agro de Limarí"; C="Coquimbo"; D=45239; E=4; F="Fruta"; G=100107; H="Otros"; I=100107002; J="Chirimoya"; K="Cultivar IV Región"; L="Segunda";  M=200; N=8000;  O=9000;  P=8500;  Q="$/bandeja 10 kilos"; R="Provincia de Limarí"; S=850;  T=10 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 4).NumberFormat = $dateFormat
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
